$wb = $excel.ActiveWorkbook

# Sheet "Overview": G2 holds the "Latest HO Xliff Generate Date" for the first file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 03:07:16"

# Sheet "zh-cn": H2 = Correspond Handoff Datetime, K2 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 03:07:12"
$wsZhCn.Range("K2").Value = "2016-08-17 03:07:28"

# Sheet "de-de": K2 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-17 03:07:35"
